$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was collected for Albahaca and needs to be
# inserted as row 420 (pushing the existing rows 420-428 down to 421-429).
$ws.Rows("420:420").Insert()

$ws.Range("A420").Value = 9
$ws.Range("B420").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C420").Value = "Metropolitana"
$ws.Range("D420").Value = 44890
$ws.Range("E420").Value = 13
$ws.Range("F420").Value = 100112052
$ws.Range("G420").Value = "Albahaca"
$ws.Range("H420").Value = "Sin especificar"
$ws.Range("I420").Value = "Primera"
$ws.Range("J420").Value = 215
$ws.Range("K420").Value = 6000
$ws.Range("L420").Value = 7000
$ws.Range("M420").Value = 6395
$ws.Range("N420").Value = "`$/docena de matas"
$ws.Range("O420").Value = "Provincia de Chacabuco"
$ws.Range("P420").Value = 1066
$ws.Range("Q420").Value = 6
$ws.Range("R420").Value = "Hortaliza"
